# Update specific numeric values in Sheet1 as per the commit "Update Name of Algo"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 5.745499999999998
$ws.Range("A3").Value = -21.28910000000003
$ws.Range("B5").Value = 4.729400000000005
$ws.Range("C5").Value = -14.34740000000001
$ws.Range("C9").Value = -12.07020000000002
$ws.Range("C11").Value = -13.3423
$ws.Range("A14").Value = -20.27449999999998
$ws.Range("A16").Value = -21.20250000000003
$ws.Range("B16").Value = 5.323200000000003
$ws.Range("C17").Value = -11.5807
$ws.Range("A21").Value = -21.25220000000001
$ws.Range("C21").Value = -11.2045
$ws.Range("A23").Value = -21.18590000000001
$ws.Range("A25").Value = -22.43030000000003
